$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the original "Meta description: ..." paragraph that currently
#    sits right under the H1 title at the top of the document.
# ---------------------------------------------------------------------------
$metaFinder = $d.Content
$foundMeta = $metaFinder.Find.Execute("Meta description", $false, $false, $false, $false, $false, `
                                       $true, 1, $false, "", 0)
if ($foundMeta) {
    $metaPar = $metaFinder.Paragraphs(1)
    $metaPar.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) The final paragraph of the document (the "DALLE, please create a
#    feature image ..." image-prompt paragraph) gets turned into two
#    paragraphs:
#       - a new bold "Play Deep Sea Magic Slot Game for Free - Review" line
#       - the (formerly "Meta description") italic blurb, now holding just
#         the descriptive sentence without the "Meta description" label
#    We replace the whole paragraph range in one shot via InsertXML so the
#    exact run layout (leading empty run + formatted run) is reproduced.
# ---------------------------------------------------------------------------
$dalleFinder = $d.Content
$null = $dalleFinder.Find.Execute("DALLE, please create a feature image", $false, $false, $false, `
                                   $false, $false, $true, 1, $false, "", 0)
$dallePar = $dalleFinder.Paragraphs(1)
$wholeRange = $dallePar.Range

$xmlFrag = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:r/>' + `
        '<w:r><w:rPr><w:b/></w:rPr><w:t>Play Deep Sea Magic Slot Game for Free - Review</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:r/>' + `
        '<w:r><w:rPr><w:i/></w:rPr><w:t>Find out about the features of the Deep Sea Magic slot game and play it for free. Read our review before playing for real money.</w:t></w:r>' + `
    '</w:p>' + `
'</pkg:xmlData>'

$null = $wholeRange.InsertXML($xmlFrag)
